$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 151 (shifts old rows 151..236 down to 152..237)
$ws.Rows(151).Insert()

# Fill in the new row 151 with data (same as the former row 151 except
# Fecha (D) and Volumen (J) which carry new values for this entry)
$ws.Range("A151").Value = 10
$ws.Range("B151").Value = "Vega Modelo de Temuco"
$ws.Range("C151").Value = "La Araucanía"
$ws.Range("D151").Value = 44606
$ws.Range("E151").Value = 9
$ws.Range("F151").Value = 100112001
$ws.Range("G151").Value = "Berenjena"
$ws.Range("H151").Value = "Sin especificar"
$ws.Range("I151").Value = "Primera"
$ws.Range("J151").Value = 100
$ws.Range("K151").Value = 10000
$ws.Range("L151").Value = 10000
$ws.Range("M151").Value = 10000
$ws.Range("N151").Value = "`$/caja 60 unidades"
$ws.Range("O151").Value = "Región del Maule"
$ws.Range("P151").Value = 167
$ws.Range("Q151").Value = 60
$ws.Range("R151").Value = "Hortaliza"
